$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# The 'Instructions' sheet is protected; temporarily unprotect so we can edit
# the locked cells, then re-protect afterwards to restore the sheet's state.
$ws.Unprotect()

# Insert a new row at position 6 (pushes old row 6 "Do not edit the other sheets." down to row 7,
# and everything below shifts down by one as well).
$ws.Rows.Item(6).Insert()

# Update version number in A2
$ws.Range("A2").Value = "Version 1.2.3"

# Update text of A5 (new instruction about consecutive rows)
$ws.Range("A5").Value = "Please use consecutive rows (no blank rows)."

# Set the new row 6 text (instruction about not editing the header row)
$ws.Range("A6").Value = "Do not edit the header row of the 'Antibodies' sheet."

# Restore sheet protection
$ws.Protect()

